# Refresh the "cryptos" price/volume table (rows 2-51) to the latest
# snapshot. Columns B/C (coin name / link) only change for the two rows
# whose rank flipped (15 <-> 16); columns D (price) and E (1h volume %)
# are refreshed throughout.
#
# Note: several Price values (column D) are plain decimals that Excel's
# input parser would normally auto-convert to a Number (e.g. "0.999",
# "586.94"), unlike the thousand-separated prices (e.g. "65.464.34")
# which stay text on their own. To keep those columns as text (matching
# the sheet's existing inlineStr cells) we enter them with a leading
# apostrophe (forces literal text, same as typing '0.999 in Excel) and
# then restore the cell's original look by copying the style from a
# neighboring, never-reformatted default-style cell (D13) so the
# quote-prefix formatting doesn't visibly change the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$styleRef = $ws.Range("D13")  # stable default-style reference cell

$ws.Range("D2").Value = "65.464.34"
$ws.Range("E2").Value = "  -0.29%  "

$ws.Range("D3").Value = "3.423.86"
$ws.Range("E3").Value = "  -2.33%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = $styleRef.Style
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'586.94"
$ws.Range("D5").Style = $styleRef.Style
$ws.Range("E5").Value = "  -1.45%  "

$ws.Range("D6").Value = "'137.72"
$ws.Range("D6").Style = $styleRef.Style
$ws.Range("E6").Value = "  -3.69%  "

$ws.Range("D7").Value = "3.422.41"
$ws.Range("E7").Value = "  -2.35%  "

$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = $styleRef.Style
$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").Value = "'0.499"
$ws.Range("D9").Style = $styleRef.Style
$ws.Range("E9").Value = "  -0.34%  "

$ws.Range("E10").Value = "  -4.92%  "

$ws.Range("E11").Value = "  -8.80%  "

$ws.Range("E12").Value = "  -7.00%  "

$ws.Range("E13").Value = "  -2.21%  "

$ws.Range("E14").Value = "  -9.61%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "3.432.97"
$ws.Range("E15").Value = "  -2.06%  "

$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "'26.25"
$ws.Range("D16").Style = $styleRef.Style
$ws.Range("E16").Value = "  -8.33%  "

$ws.Range("D17").Value = "65.408.22"
$ws.Range("E17").Value = "  -0.42%  "

$ws.Range("E18").Value = "  -1.63%  "

$ws.Range("D19").Value = "'9.81"
$ws.Range("D19").Style = $styleRef.Style
$ws.Range("E19").Value = "  -10.11%  "

$ws.Range("E20").Value = "  -4.96%  "

$ws.Range("E21").Value = "  -4.86%  "

$ws.Range("D22").Value = "'391.21"
$ws.Range("D22").Style = $styleRef.Style
$ws.Range("E22").Value = "  -5.08%  "

$ws.Range("E23").Value = "  -6.64%  "

$ws.Range("D24").Value = "'73.16"
$ws.Range("D24").Style = $styleRef.Style
$ws.Range("E24").Value = "  -5.37%  "

$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("D26").Value = "3.559.39"
$ws.Range("E26").Value = "  -2.11%  "

$ws.Range("E27").Value = "  -7.22%  "

$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("E29").Value = "  -6.66%  "

$ws.Range("E30").Value = "  -9.17%  "

$ws.Range("E31").Value = "  -8.36%  "

$ws.Range("D32").Value = "3.428.73"
$ws.Range("E32").Value = "  -1.99%  "

$ws.Range("E34").Value = "  -5.51%  "

$ws.Range("D35").Value = "'22.95"
$ws.Range("D35").Style = $styleRef.Style
$ws.Range("E35").Value = "  -5.20%  "

$ws.Range("D36").Value = "'172.63"
$ws.Range("D36").Style = $styleRef.Style
$ws.Range("E36").Value = "  -1.46%  "

$ws.Range("D37").Value = "'6.85"
$ws.Range("D37").Style = $styleRef.Style
$ws.Range("E37").Value = "  -8.18%  "

$ws.Range("E38").Value = "  -7.05%  "

$ws.Range("D39").Value = "'1.47"
$ws.Range("D39").Style = $styleRef.Style
$ws.Range("E39").Value = "  -6.57%  "

$ws.Range("E40").Value = "  -8.04%  "

$ws.Range("D41").Value = "'0.0763"
$ws.Range("D41").Style = $styleRef.Style
$ws.Range("E41").Value = "  -6.44%  "

$ws.Range("D42").Value = "'0.819"
$ws.Range("D42").Style = $styleRef.Style
$ws.Range("E42").Value = "  -3.75%  "

$ws.Range("D43").Value = "'43.61"
$ws.Range("D43").Style = $styleRef.Style
$ws.Range("E43").Value = "  -3.45%  "

$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("D45").Value = "'4.42"
$ws.Range("D45").Style = $styleRef.Style
$ws.Range("E45").Value = "  -12.02%  "

$ws.Range("D46").Value = "'1.61"
$ws.Range("D46").Style = $styleRef.Style
$ws.Range("E46").Value = "  -8.97%  "

$ws.Range("E47").Value = "  +3.66%  "

$ws.Range("D48").Value = "'22.32"
$ws.Range("D48").Style = $styleRef.Style
$ws.Range("E48").Value = "  +0.51%  "

$ws.Range("D49").Value = "'6.52"
$ws.Range("D49").Style = $styleRef.Style
$ws.Range("E49").Value = "  -7.64%  "

$ws.Range("D50").Value = "'2.09"
$ws.Range("D50").Style = $styleRef.Style
$ws.Range("E50").Value = "  -12.81%  "

$ws.Range("D51").Value = "2.192.70"
$ws.Range("E51").Value = "  -6.69%  "
